$d = $word.ActiveDocument

$replacements = @(
    @{ old = "57×79="; new = "49×15=" },
    @{ old = "73×50="; new = "15×40=" },
    @{ old = "42×81="; new = "83×16=" },
    @{ old = "85×12="; new = "37×60=" },
    @{ old = "24×24="; new = "48×58=" },
    @{ old = "81×61="; new = "46×39=" },
    @{ old = "24×46="; new = "85×92=" },
    @{ old = "48×34="; new = "28×37=" },
    @{ old = "38×29="; new = "70×48=" },
    @{ old = "57×87="; new = "67×34=" },
    @{ old = "22×72="; new = "53×68=" },
    @{ old = "60×51="; new = "83×57=" },
    @{ old = "70×79="; new = "12×52=" },
    @{ old = "53×24="; new = "41×77=" },
    @{ old = "55×73="; new = "89×79=" },
    @{ old = "80×51="; new = "19×37=" },
    @{ old = "77×67="; new = "59×55=" },
    @{ old = "61×95="; new = "13×40=" },
    @{ old = "72×11="; new = "45×81=" },
    @{ old = "57×18="; new = "64×74=" },
    @{ old = "56×55="; new = "52×22=" },
    @{ old = "28×86="; new = "54×31=" },
    @{ old = "92×69="; new = "55×43=" },
    @{ old = "41×22="; new = "53×58=" },
    @{ old = "99×92="; new = "85×83=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
